# Staging.Activity.xlsx - header columns were re-mapped and one field was
# renamed (ProjectSourceKey -> ProjectBusinessKey) as the staging templates
# moved to their new directory.
#
# Row 2 holds the column headers used by the staging importer. Column B
# (previously "ShortName") now carries "Code", column D (previously
# "TextDescription") now carries the renamed "ProjectBusinessKey" field,
# and the old B/D values slide over to E/F respectively.
#
# Resulting header row, left to right:
#   A2 Activity_ID  (unchanged)
#   B2 Code
#   C2 LongName     (unchanged)
#   D2 ProjectBusinessKey
#   E2 ShortName
#   F2 TextDescription

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = "Code"
$ws.Range("D2").Value2 = "ProjectBusinessKey"
$ws.Range("E2").Value2 = "ShortName"
$ws.Range("F2").Value2 = "TextDescription"

# The author's selection when they saved the file moved from C10 to D16.
$ws.Range("D16").Select()
